$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename a handful of item labels (units/spelling tweaks) ---
$ws.Range("A8").Value  = "P--Kreatinin (enz) (mikromol/L)"
$ws.Range("A16").Value = "B--Leukocyter (x10(9)/L)"
$ws.Range("A17").Value = "B--Erytrocyter (x10(12)/L)"
$ws.Range("A19").Value = "B--EVF"
$ws.Range("A22").Value = "P--Glukos (mmol/L)"

# --- 2. P--Natrium (row 21) previously had no Min/Max values - fill them in ---
$ws.Range("B21").Value = 137
$ws.Range("C21").Value = 145

# --- 3. Append four new reference rows ---
$ws.Range("A24").Value = "B--Trombocyter (x10(9)/L)"
$ws.Range("B24").Value = 145
$ws.Range("C24").Value = 348

$ws.Range("A25").Value = "B--Neutrofila granulocyter (x10(9)/L)"
$ws.Range("B25").Value = 1.6
$ws.Range("C25").Value = 5.9

$ws.Range("A26").Value = "B--Myelocyter (x10(9)/L)"
$ws.Range("B26").Value = -999
$ws.Range("C26").Value = 0

$ws.Range("A27").Value = "B--Metamyelocyter (x10(9)/L)"
$ws.Range("B27").Value = -999
$ws.Range("C27").Value = 0

# --- 4. Add a bold, boxed divider style under row 16 (A16:D16) ---
$divider = $ws.Range("A16:D16")
$divider.Font.Bold = $true

# Accent 5 theme blue (5B9BD5) expressed as a BGR long for the legacy .Color setter
$accentBlue = 13998939

$left = $divider.Borders.Item(7)
$left.LineStyle = 1
$left.Weight = 2
$left.Color = $accentBlue

$right = $divider.Borders.Item(10)
$right.LineStyle = 1
$right.Weight = 2
$right.Color = $accentBlue

$top = $divider.Borders.Item(8)
$top.LineStyle = 1
$top.Weight = 2
$top.Color = $accentBlue

$bottom = $divider.Borders.Item(9)
$bottom.LineStyle = 1
$bottom.Weight = -4138
$bottom.Color = $accentBlue

# --- 5. Move the active selection to the last entered row ---
$ws.Range("A26").Select()
